# Generate Report for Handoff
#
# The source file that was previously pending handoff
# (eb02f3a4-3b8c-4b21-8d75-c1fff3ed7e3b.md) has been renamed/re-handed-off as
# 0db4b02f-32a6-4080-8f9f-657f4213dcc2.md (new handoff package hash
# 37b5125863226766a3cd29cc09197671830c6786, refreshed handoff timestamps),
# and a brand-new file ffffa99612f2-941e-44be-9143-26876ac55760.md has
# entered the pipeline with status "Ready for handoff" using that same
# handoff package. This inserts one new row into every sheet for the new
# file (pushing the ".localization-config" bookkeeping row down by one).

$wb = $excel.ActiveWorkbook

$oldFile = "eb02f3a4-3b8c-4b21-8d75-c1fff3ed7e3b.md"
$newFile = "0db4b02f-32a6-4080-8f9f-657f4213dcc2.md"
$addedFile = "ffffa99612f2-941e-44be-9143-26876ac55760.md"
$configFile = ".localization-config"

$oldZhXlf = "eb02f3a4-3b8c-4b21-8d75-c1fff3ed7e3b.d5209ffb99d0038c24af071d96c7e9a21945fa85.zh-cn.xlf"
$newZhXlf = "0db4b02f-32a6-4080-8f9f-657f4213dcc2.37b5125863226766a3cd29cc09197671830c6786.zh-cn.xlf"
$oldDeXlf = "eb02f3a4-3b8c-4b21-8d75-c1fff3ed7e3b.d5209ffb99d0038c24af071d96c7e9a21945fa85.de-de.xlf"
$newDeXlf = "0db4b02f-32a6-4080-8f9f-657f4213dcc2.37b5125863226766a3cd29cc09197671830c6786.de-de.xlf"

$zhHandoffTime = "2016-02-29 13:38:47"
$deHandoffTime = "2016-02-29 13:38:58"
$epoch = "0001-01-01 00:00:00"

$blueUnderline = 15570276   # FF6495ED, expressed BGR for Font.Color

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $blueUnderline
}

function Base-Url($repoPath) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/d1e363b25f0c08df0c6d8992782b3333cd3dc00c/" + $repoPath
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows(3).Insert()

$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = $addedFile
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("A4").Value = $configFile
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), (Base-Url "e2e/$newFile"), "", "", $newFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), (Base-Url "e2e/$addedFile"), "", "", $addedFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), (Base-Url $configFile), "", "", $configFile)

Style-AsHyperlink $wsOverview.Range("A2")
Style-AsHyperlink $wsOverview.Range("A3")
Style-AsHyperlink $wsOverview.Range("A4")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows(3).Insert()

$wsZh.Range("A2").Value = $newFile
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $zhHandoffTime
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = $addedFile
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $newZhXlf
$wsZh.Range("D3").Value = $zhHandoffTime
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("A4").Value = $configFile
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), (Base-Url "e2e/$newFile"), "", "", $newFile)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d06f56147ff8ad447fdc45a1c373cbe69a2fc7b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZhXlf", "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), (Base-Url "e2e/$addedFile"), "", "", $addedFile)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d06f56147ff8ad447fdc45a1c373cbe69a2fc7b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZhXlf", "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), (Base-Url $configFile), "", "", $configFile)

Style-AsHyperlink $wsZh.Range("A2")
Style-AsHyperlink $wsZh.Range("C2")
Style-AsHyperlink $wsZh.Range("A3")
Style-AsHyperlink $wsZh.Range("C3")
Style-AsHyperlink $wsZh.Range("A4")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows(3).Insert()

$wsDe.Range("A2").Value = $newFile
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $deHandoffTime
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = $addedFile
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $newDeXlf
$wsDe.Range("D3").Value = $deHandoffTime
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("A4").Value = $configFile
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), (Base-Url "e2e/$newFile"), "", "", $newFile)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd5c7cdab67f4093681f72fe306271bfe6b3411e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDeXlf", "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), (Base-Url "e2e/$addedFile"), "", "", $addedFile)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd5c7cdab67f4093681f72fe306271bfe6b3411e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDeXlf", "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), (Base-Url $configFile), "", "", $configFile)

Style-AsHyperlink $wsDe.Range("A2")
Style-AsHyperlink $wsDe.Range("C2")
Style-AsHyperlink $wsDe.Range("A3")
Style-AsHyperlink $wsDe.Range("C3")
Style-AsHyperlink $wsDe.Range("A4")
